$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = 1.65
$ws.Range("H5").Value = 3.5
$ws.Range("I5").Value = 5.75
$ws.Range("J5").Value = 2.4
$ws.Range("K5").Value = 1.83
$ws.Range("L5").Value = 8
$ws.Range("M5").Value = 1.16
$ws.Range("N5").Value = 4.65
$ws.Range("O5").Value = 1.7
$ws.Range("P5").Value = 2.05
$ws.Range("W5").Value = 5.6
$ws.Range("X5").Value = 1.11
$ws.Range("AA5").Value = 2.82
$ws.Range("AB5").Value = 1.38
$ws.Range("AC5").Value = 4.33
$ws.Range("AD5").Value = 6
$ws.Range("AF5").Value = 12
$ws.Range("AJ5").Value = 8
$ws.Range("AK5").Value = 34
$ws.Range("AL5").Value = 151
$ws.Range("AN5").Value = 9
$ws.Range("AP5").Value = 21
$ws.Range("AQ5").Value = 81
$ws.Range("AR5").Value = 67
$ws.Range("AS5").Value = 81
$ws.Range("G7").Value = 1.72
$ws.Range("H7").Value = 3.15
$ws.Range("I7").Value = 5.6
$ws.Range("J7").Value = 2.27
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 5.9
$ws.Range("M7").Value = 1.11
$ws.Range("N7").Value = 5.6
$ws.Range("O7").Value = 1.47
$ws.Range("P7").Value = 2.52
$ws.Range("S7").Value = 2.35
$ws.Range("T7").Value = 1.53
$ws.Range("W7").Value = 4.1
$ws.Range("X7").Value = 1.19
$ws.Range("Y7").Value = 1.5
$ws.Range("Z7").Value = 2.42
$ws.Range("AA7").Value = 2.15
$ws.Range("AB7").Value = 1.62
$ws.Range("AC7").Value = 5.2
$ws.Range("AE7").Value = 8.5
$ws.Range("AF7").Value = 13.5
$ws.Range("AG7").Value = 16
$ws.Range("AH7").Value = 37
$ws.Range("AI7").Value = 5.6
$ws.Range("AJ7").Value = 6.3
$ws.Range("AK7").Value = 20
$ws.Range("AL7").Value = 120
$ws.Range("AM7").Value = 900
$ws.Range("AN7").Value = 10.75
$ws.Range("AP7").Value = 19
$ws.Range("AR7").Value = 80
$ws.Range("AS7").Value = 90
$ws.Range("G8").Value = 3
$ws.Range("H8").Value = 2.55
$ws.Range("I8").Value = 2.85
$ws.Range("J8").Value = 3.9
$ws.Range("L8").Value = 3.7
$ws.Range("M8").Value = 1.19
$ws.Range("N8").Value = 4.2
$ws.Range("O8").Value = 1.78
$ws.Range("P8").Value = 1.93
$ws.Range("S8").Value = 3.25
$ws.Range("T8").Value = 1.3
$ws.Range("W8").Value = 6
$ws.Range("X8").Value = 1.09
$ws.Range("Y8").Value = 1.75
$ws.Range("Z8").Value = 1.98
$ws.Range("AC8").Value = 5.7
$ws.Range("AD8").Value = 13
$ws.Range("AE8").Value = 12.5
$ws.Range("AF8").Value = 40
$ws.Range("AG8").Value = 40
$ws.Range("AH8").Value = 75
$ws.Range("AI8").Value = 4.2
$ws.Range("AJ8").Value = 5.5
$ws.Range("AK8").Value = 22
$ws.Range("AN8").Value = 5.6
$ws.Range("AO8").Value = 12
$ws.Range("AP8").Value = 12
$ws.Range("AQ8").Value = 37
$ws.Range("AR8").Value = 37
$ws.Range("AS8").Value = 70
$ws.Range("G9").Value = 1.38
$ws.Range("H9").Value = 4.15
$ws.Range("I9").Value = 8.75
$ws.Range("J9").Value = 1.91
$ws.Range("K9").Value = 2.2
$ws.Range("L9").Value = 7.4
$ws.Range("M9").Value = 1.05
$ws.Range("N9").Value = 7.7
$ws.Range("O9").Value = 1.26
$ws.Range("P9").Value = 3.5
$ws.Range("S9").Value = 1.78
$ws.Range("T9").Value = 1.93
$ws.Range("W9").Value = 2.85
$ws.Range("X9").Value = 1.37
$ws.Range("Y9").Value = 1.4
$ws.Range("Z9").Value = 2.75
$ws.Range("AA9").Value = 2.02
$ws.Range("AB9").Value = 1.7
$ws.Range("AC9").Value = 6.1
$ws.Range("AD9").Value = 6.1
$ws.Range("AE9").Value = 8.25
$ws.Range("AG9").Value = 11.75
$ws.Range("AH9").Value = 29
$ws.Range("AI9").Value = 7.7
$ws.Range("AJ9").Value = 8.25
$ws.Range("AK9").Value = 19.5
$ws.Range("AL9").Value = 100
$ws.Range("AM9").Value = 800
$ws.Range("AN9").Value = 22
$ws.Range("AP9").Value = 26
$ws.Range("AQ9").Value = 250
$ws.Range("AR9").Value = 100
$ws.Range("AS9").Value = 80
$ws.Range("G12").Value = 1.2
$ws.Range("H12").Value = 6.25
$ws.Range("J12").Value = 1.57
$ws.Range("L12").Value = 9
$ws.Range("S12").Value = 1.44
$ws.Range("T12").Value = 2.63
$ws.Range("W12").Value = 2.1
$ws.Range("X12").Value = 1.67
$ws.Range("AA12").Value = 2.1
$ws.Range("AB12").Value = 1.67
$ws.Range("AC12").Value = 8.5
$ws.Range("AJ12").Value = 13
$ws.Range("AK12").Value = 26
$ws.Range("AM12").Value = 900
$ws.Range("AN12").Value = 29
$ws.Range("S13").Value = 1.53
$ws.Range("T13").Value = 2.38
$ws.Range("U13").Value = 1.88
$ws.Range("V13").Value = 1.93
$ws.Range("G15").Value = 1.67
$ws.Range("H15").Value = 4.2
$ws.Range("J15").Value = 2.2
$ws.Range("K15").Value = 2.5
$ws.Range("L15").Value = 4.5
$ws.Range("N15").Value = 17
$ws.Range("O15").Value = 1.17
$ws.Range("P15").Value = 5
$ws.Range("S15").Value = 1.53
$ws.Range("T15").Value = 2.4
$ws.Range("U15").Value = 1.88
$ws.Range("V15").Value = 1.98
$ws.Range("W15").Value = 2.25
$ws.Range("X15").Value = 1.57
$ws.Range("Y15").Value = 1.29
$ws.Range("Z15").Value = 3.5
$ws.Range("AA15").Value = 1.57
$ws.Range("AB15").Value = 2.25
$ws.Range("AC15").Value = 10
$ws.Range("AD15").Value = 9.5
$ws.Range("AG15").Value = 12
$ws.Range("AI15").Value = 17
$ws.Range("AJ15").Value = 8.5
$ws.Range("AM15").Value = 126
$ws.Range("AN15").Value = 17
$ws.Range("AR15").Value = 29
$ws.Range("AS15").Value = 29
$ws.Range("G17").Value = 5.7
$ws.Range("H17").Value = 4.2
$ws.Range("I17").Value = 1.5
$ws.Range("J17").Value = 5.4
$ws.Range("K17").Value = 2.37
$ws.Range("Y17").Value = 1.32
$ws.Range("Z17").Value = 3.1
$ws.Range("AC17").Value = 17
$ws.Range("AD17").Value = 35
$ws.Range("AF17").Value = 110
$ws.Range("AH17").Value = 50
$ws.Range("AN17").Value = 7.7
$ws.Range("AO17").Value = 7.5
$ws.Range("AQ17").Value = 10.5
$ws.Range("AR17").Value = 11.5
$ws.Range("AS17").Value = 23
$ws.Range("G18").Value = 1.82
$ws.Range("H18").Value = 3.6
$ws.Range("I18").Value = 3.85
$ws.Range("J18").Value = 2.4
$ws.Range("K18").Value = 2.15
$ws.Range("L18").Value = 4.25
$ws.Range("O18").Value = 1.28
$ws.Range("P18").Value = 3.05
$ws.Range("S18").Value = 1.82
$ws.Range("T18").Value = 1.78
$ws.Range("W18").Value = 2.9
$ws.Range("X18").Value = 1.31
$ws.Range("Y18").Value = 1.38
$ws.Range("Z18").Value = 2.6
$ws.Range("AA18").Value = 1.75
$ws.Range("AB18").Value = 1.85
$ws.Range("AC18").Value = 7.2
$ws.Range("AD18").Value = 8.5
$ws.Range("AE18").Value = 8.25
$ws.Range("AF18").Value = 15
$ws.Range("AG18").Value = 14.5
$ws.Range("AH18").Value = 27
$ws.Range("AI18").Value = 10.5
$ws.Range("AJ18").Value = 6.9
$ws.Range("AK18").Value = 15.5
$ws.Range("AL18").Value = 70
$ws.Range("AM18").Value = 600
$ws.Range("AN18").Value = 11
$ws.Range("AO18").Value = 21
$ws.Range("AP18").Value = 13
$ws.Range("AQ18").Value = 55
$ws.Range("AR18").Value = 35
